# Insert a new weekly record at row 575, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 575. This shifts rows 575..633 down to 576..634
# and Excel automatically preserves per-column formatting (e.g. the date style on column D).
$ws.Rows(575).Insert()

# Populate the newly inserted row 575 with the new weekly record.
$ws.Range("A575").Value = 5
$ws.Range("B575").Value = "Macroferia Regional de Talca"
$ws.Range("C575").Value = "Maule"
$ws.Range("D575").Value = 45194
$ws.Range("E575").Value = 7
$ws.Range("F575").Value = 100112023
$ws.Range("G575").Value = "Brócoli"
$ws.Range("H575").Value = "Sin especificar"
$ws.Range("I575").Value = "Primera"
$ws.Range("J575").Value = 3000
$ws.Range("K575").Value = 1200
$ws.Range("L575").Value = 1200
$ws.Range("M575").Value = 1200
$ws.Range("N575").Value = "$/unidad"
$ws.Range("O575").Value = "Región del Maule"
$ws.Range("P575").Value = 1200
$ws.Range("Q575").Value = 1
$ws.Range("R575").Value = "Hortaliza"
